{"js": "// Applies the textual edits from the commit:\n//  - \"huisstyle\" -> \"huisstijl\"\n//  - \"schitterende design\" -> \"schitterend design\"\n//  - \"Jp\" -> \"JP\" (all whole-word occurrences)\n//  - \"Verbeteringen voor de opdrachtgever\" -> \"Verbeteringen voor de toekomst\"\n//    (both the heading and its Table-of-Contents entry)\n//  - final paragraph body text fully rewritten\nconst body = context.document.body;\n\nasync function replaceAll(find, repl, opts) {\n  const results = body.search(find, opts || { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(repl, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Spelling fixes in \"Reflectie over het product\".\nawait replaceAll(\"huisstyle\", \"huisstijl\");\nawait replaceAll(\"schitterende design\", \"schitterend design\");\n\n// \"Jp\" -> \"JP\" everywhere it appears as a whole word (4x, in \"Reflectie over proces\").\nawait replaceAll(\"Jp\", \"JP\", { matchCase: true, matchWholeWord: true });\n\n// Heading rename (also updates the matching Table of Contents entry text).\nawait replaceAll(\n  \"Verbeteringen voor de opdrachtgever\",\n  \"Verbeteringen voor de toekomst\"\n);\n\n// Final paragraph (\"Verbeteringen voor de toekomst\" section body) rewritten entirely.\nawait replaceAll(\n  \"Emails op tijd antwoorden en meer docent om zowel technische vragen als niet technische vragen te beantwoorden.\",\n  \"Om eerlijk te zijn kan ik mezelf niet verbeteren, ik heb optimaal gepresteerd in alle factoren. \" +\n  \"Zowel documentatie als in het maken van de applicatie. In kon gewoon kan gewoon niet beter zijn, \" +\n  \"ik ben altijd eerlijk, probeer het maximale uit mensen te halen en uit mijn werk. Ik ben nog nooit \" +\n  \"zo hoogmoedig om dit te zeggen. Ik kan alleen mijn kennis vergroten en dat zal ik doen in de \" +\n  \"aankomende periodes.\"\n);\n", "ps1": "# Applies the textual edits from the commit:\n#  - \"huisstyle\" -> \"huisstijl\"\n#  - \"schitterende design\" -> \"schitterend design\"\n#  - \"Jp\" -> \"JP\" (all whole-word occurrences)\n#  - \"Verbeteringen voor de opdrachtgever\" -> \"Verbeteringen voor de toekomst\"\n#    (both the heading and its Table-of-Contents entry)\n#  - final paragraph body text fully rewritten\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText, $matchWholeWord) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = [bool]$matchWholeWord\n    $range.Find.Wrap = 1            # wdFindContinue\n    $range.Find.Forward = $true\n    $range.Find.Format = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Spelling fixes in \"Reflectie over het product\".\nReplace-All \"huisstyle\" \"huisstijl\" $false\nReplace-All \"schitterende design\" \"schitterend design\" $false\n\n# \"Jp\" -> \"JP\" everywhere it appears as a whole word (4x, in \"Reflectie over proces\").\nReplace-All \"Jp\" \"JP\" $true\n\n# Heading rename (also updates the matching Table of Contents entry text).\nReplace-All \"Verbeteringen voor de opdrachtgever\" \"Verbeteringen voor de toekomst\" $false\n\n# Final paragraph (\"Verbeteringen voor de toekomst\" section body) rewritten entirely.\n$newClosing = \"Om eerlijk te zijn kan ik mezelf niet verbeteren, ik heb optimaal gepresteerd in alle factoren. \" + `\n  \"Zowel documentatie als in het maken van de applicatie. In kon gewoon kan gewoon niet beter zijn, \" + `\n  \"ik ben altijd eerlijk, probeer het maximale uit mensen te halen en uit mijn werk. Ik ben nog nooit \" + `\n  \"zo hoogmoedig om dit te zeggen. Ik kan alleen mijn kennis vergroten en dat zal ik doen in de \" + `\n  \"aankomende periodes.\"\nReplace-All \"Emails op tijd antwoorden en meer docent om zowel technische vragen als niet technische vragen te beantwoorden.\" $newClosing $false\n"}
